$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.434.29"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "3.281.37"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'253.69"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'622.33"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'1.43"
$ws.Range("E7").Value = "  +21.83%  "
$ws.Range("D8").Value = "'0.400"
$ws.Range("E8").Value = "  +3.88%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +21.44%  "
$ws.Range("D11").Value = "3.281.24"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "'39.43"
$ws.Range("E13").Value = "  +10.32%  "
$ws.Range("D14").Value = "99.172.27"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "'0.0000247"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "3.869.89"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "'5.47"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "3.277.70"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("E19").Value = "  -4.80%  "
$ws.Range("D20").Value = "'15.24"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("E21").Value = "  +8.42%  "
$ws.Range("D22").Value = "'486.50"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "'5.61"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'89.05"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "'0.327"
$ws.Range("E27").Value = "  +34.07%  "
$ws.Range("D28").Value = "'11.95"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "3.424.74"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("D32").Value = "'0.136"
$ws.Range("E32").Value = "  +10.44%  "
$ws.Range("D33").Value = "'10.33"
$ws.Range("E33").Value = "  +12.06%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "'27.87"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").Value = "'0.149"
$ws.Range("D37").Value = "'0.473"
$ws.Range("E37").Value = "  +5.98%  "
$ws.Range("D38").Value = "'7.18"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'24.82"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  -5.54%  "
$ws.Range("D42").Value = "'3.69"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'0.770"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "'3.08"
$ws.Range("E46").Value = "  -5.61%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").Value = "'157.24"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("E49").Value = "  +7.84%  "
$ws.Range("E50").Value = "  +15.01%  "
$ws.Range("D51").Value = "'4.70"
$ws.Range("E51").Value = "  +4.52%  "
